# Applies the edit described in the commit: updates two odds values in row 3
# (Talleres Cordoba vs Lanus) and inserts a new match row (Platense vs Dep.
# Riestra) as row 5, pushing the previous rows 5-7 down to rows 6-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Talleres Cordoba vs Lanus): update Odd_Over15_FT / Odd_Under15_FT ---
$ws.Cells.Item(3, 15).Value = 1.5   # O3 Odd_Over15_FT
$ws.Cells.Item(3, 16).Value = 2.5   # P3 Odd_Under15_FT

# --- Insert a new row 5, shifting current rows 5,6,7 down to 6,7,8 ---
$ws.Rows("5:5").Insert()

# --- Fill new row 5: ARGENTINA - TORNEO BETANO, Platense vs Dep. Riestra ---
$ws.Cells.Item(5, 1).Value = "AgfqVj1p"       # Id
$ws.Cells.Item(5, 2).NumberFormat = "@"       # Date
$ws.Cells.Item(5, 2).Value = "07/11/2024"
$ws.Cells.Item(5, 2).ClearFormats()
$ws.Cells.Item(5, 3).Value = "21:00"       # Time
$ws.Cells.Item(5, 4).Value = "ARGENTINA - TORNEO BETANO"       # League
$ws.Cells.Item(5, 5).Value = "Platense"       # Home
$ws.Cells.Item(5, 6).Value = "Dep. Riestra"       # Away
$ws.Cells.Item(5, 7).Value = 1.8       # Odd_H_FT
$ws.Cells.Item(5, 8).Value = 3.1       # Odd_D_FT
$ws.Cells.Item(5, 9).Value = 5.75       # Odd_A_FT
$ws.Cells.Item(5, 10).Value = 2.6       # Odd_H_HT
$ws.Cells.Item(5, 11).Value = 1.83       # Odd_D_HT
$ws.Cells.Item(5, 12).Value = 6.5       # Odd_A_HT
$ws.Cells.Item(5, 13).Value = 1.17       # Odd_Over05_FT
$ws.Cells.Item(5, 14).Value = 5       # Odd_Under05_FT
$ws.Cells.Item(5, 15).Value = 1.62       # Odd_Over15_FT
$ws.Cells.Item(5, 16).Value = 2.2       # Odd_Under15_FT
$ws.Cells.Item(5, 17).Value = 3.1       # Odd_Over25_FT
$ws.Cells.Item(5, 18).Value = 1.36       # Odd_Under25_FT
$ws.Cells.Item(5, 19).Value = 1.67       # Odd_Over05_HT
$ws.Cells.Item(5, 20).Value = 2.1       # Odd_Under05_HT
$ws.Cells.Item(5, 21).Value = 2.63       # Odd_BTTS_Yes
$ws.Cells.Item(5, 22).Value = 1.44       # Odd_BTTS_No
$ws.Cells.Item(5, 23).Value = 4.5       # Odd_CS_1-0
$ws.Cells.Item(5, 24).Value = 6.5       # Odd_CS_2-0
$ws.Cells.Item(5, 25).Value = 10       # Odd_CS_2-1
$ws.Cells.Item(5, 26).Value = 15       # Odd_CS_3-0
$ws.Cells.Item(5, 27).Value = 21       # Odd_CS_3-1
$ws.Cells.Item(5, 28).Value = 41       # Odd_CS_3-2
$ws.Cells.Item(5, 29).Value = 5       # Odd_CS_0-0
$ws.Cells.Item(5, 30).Value = 6.5       # Odd_CS_1-1
$ws.Cells.Item(5, 31).Value = 26       # Odd_CS_2-2
$ws.Cells.Item(5, 32).Value = 126       # Odd_CS_3-3
$ws.Cells.Item(5, 33).Value = 501       # Odd_CS_4-4
$ws.Cells.Item(5, 34).Value = 9.5       # Odd_CS_0-1
$ws.Cells.Item(5, 35).Value = 26       # Odd_CS_0-2
$ws.Cells.Item(5, 36).Value = 21       # Odd_CS_1-2
$ws.Cells.Item(5, 37).Value = 67       # Odd_CS_0-3
$ws.Cells.Item(5, 38).Value = 51       # Odd_CS_1-3
$ws.Cells.Item(5, 39).Value = 81       # Odd_CS_2-3
$ws.Cells.Item(5, 40).Value = 3.5       # Odd_CS_1-0_HT
$ws.Cells.Item(5, 41).Value = 11       # Odd_CS_2-0_HT
$ws.Cells.Item(5, 42).Value = 29       # Odd_CS_2-1_HT
$ws.Cells.Item(5, 43).Value = 41       # Odd_CS_3-0_HT
$ws.Cells.Item(5, 44).Value = 81       # Odd_CS_3-1_HT
$ws.Cells.Item(5, 45).Value = 401       # Odd_CS_3-2_HT
$ws.Cells.Item(5, 46).Value = 2.1       # Odd_CS_0-0_HT
$ws.Cells.Item(5, 47).Value = 11       # Odd_CS_1-1_HT
$ws.Cells.Item(5, 48).Value = 101       # Odd_CS_2-2_HT
$ws.Cells.Item(5, 49).Value = 7       # Odd_CS_0-1_HT
$ws.Cells.Item(5, 50).Value = 34       # Odd_CS_0-2_HT
$ws.Cells.Item(5, 51).Value = 51       # Odd_CS_1-2_HT
$ws.Cells.Item(5, 52).Value = 151       # Odd_CS_0-3_HT
$ws.Cells.Item(5, 53).Value = 251       # Odd_CS_1-3_HT
$ws.Cells.Item(5, 54).Value = 501       # Odd_CS_2-3_HT
$ws.Cells.Item(5, 55).Value = 126       # Odd_CS_3-3_HT
$ws.Cells.Item(5, 56).Value = 126       # Odd_CS_4-4_HT
